# Update latest output (run 20)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E4").Value = 495.6750825
$wsSchedule.Range("F4").Value = 29.14021649029982
$wsSchedule.Range("E5").Value = 417.760239
$wsSchedule.Range("F5").Value = 12.27984241622575

# --- Sheet "Detailed" ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B37").Value = 38.59672
$wsDetailed.Range("B38").Value = 38.15612

$wsDetailed.Range("B39").Value = 35.88
$wsDetailed.Range("C39").Value = "historical"

$wsDetailed.Range("B40").Value = 36.0601
$wsDetailed.Range("C40").Value = "historical"

$wsDetailed.Range("B41").Value = 80.02
$wsDetailed.Range("B42").Value = 85.95
$wsDetailed.Range("B43").Value = 85.95
$wsDetailed.Range("B44").Value = 85.95
$wsDetailed.Range("B45").Value = 80.02

$wsDetailed.Range("B47").Value = 56.98
$wsDetailed.Range("B48").Value = 56.98
$wsDetailed.Range("B49").Value = 56.98
$wsDetailed.Range("B50").Value = 58.17198
$wsDetailed.Range("B51").Value = 58.21072
$wsDetailed.Range("B52").Value = 57.06003
$wsDetailed.Range("B53").Value = 50.04197

$wsDetailed.Range("B56").Value = 56.98
$wsDetailed.Range("B57").Value = 51.17561
$wsDetailed.Range("B58").Value = 50.87011
$wsDetailed.Range("B59").Value = 64.56525000000001
$wsDetailed.Range("B60").Value = 64.45856000000001
$wsDetailed.Range("B61").Value = 75.84058
$wsDetailed.Range("B62").Value = 65

$wsDetailed.Range("B66").Value = 41.60172
$wsDetailed.Range("B67").Value = 36.06

$wsDetailed.Range("B73").Value = 36.06
$wsDetailed.Range("B74").Value = 6.4549
$wsDetailed.Range("B75").Value = 0.51

$wsDetailed.Range("B77").Value = -3.6481
$wsDetailed.Range("B78").Value = -4.81333
$wsDetailed.Range("B79").Value = 6.79107
$wsDetailed.Range("B80").Value = 4.96539
$wsDetailed.Range("B81").Value = 16.0108
$wsDetailed.Range("B82").Value = 17.3123
$wsDetailed.Range("B83").Value = 8.53261
$wsDetailed.Range("B84").Value = -11.45546
$wsDetailed.Range("B85").Value = -10.83089
$wsDetailed.Range("B86").Value = -6.88049
$wsDetailed.Range("B87").Value = -6
$wsDetailed.Range("B88").Value = -3.09313
$wsDetailed.Range("B89").Value = 29.85322

$wsDetailed.Range("B92").Value = 29.85322
$wsDetailed.Range("B93").Value = 78
